$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

# Add the new rows for the air attacker ally ("Spear Guy")
$ws.Range("A22").Value = "unitAllySpearman"
$ws.Range("B22").Value = "Spear Guy"
$ws.Range("A23").Value = "unitAllySpearmanDesc"
$ws.Range("B23").Value = "He spears in air."

# Update the view: scroll so row 7 is at top, and select A23 (the new last row)
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("A23").Select()
